$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_9_5_0"
$ws.Cells.Item(2, 2).Value = 0.824505913402631
$ws.Cells.Item(2, 3).Value = 0.7714446812286466
$ws.Cells.Item(2, 4).Value = 0.5056059913408425
$ws.Cells.Item(2, 5).Value = 0.6659461179673163
$ws.Cells.Item(2, 6).Value = 0.1942202299833298
$ws.Cells.Item(2, 7).Value = 0.238774299621582
$ws.Cells.Item(2, 8).Value = 0.4835389852523804
$ws.Cells.Item(2, 9).Value = 0.3539574444293976

$ws.Cells.Item(3, 1).Value = "model_9_5_20"
$ws.Cells.Item(3, 2).Value = 0.8370256563421017
$ws.Cells.Item(3, 3).Value = 0.6944868456081168
$ws.Cells.Item(3, 4).Value = 0.2843430719266707
$ws.Cells.Item(3, 5).Value = 0.5296644144321598
$ws.Cells.Item(3, 6).Value = 0.180364578962326
$ws.Cells.Item(3, 7).Value = 0.3191729784011841
$ws.Cells.Item(3, 8).Value = 0.6999437808990479
$ws.Cells.Item(3, 9).Value = 0.4983590543270111

$ws.Cells.Item(4, 1).Value = "model_9_5_19"
$ws.Cells.Item(4, 2).Value = 0.8384879696246339
$ws.Cells.Item(4, 3).Value = 0.6942746768925516
$ws.Cells.Item(4, 4).Value = 0.2965814975077447
$ws.Cells.Item(4, 5).Value = 0.534869691937284
$ws.Cells.Item(4, 6).Value = 0.178746223449707
$ws.Cells.Item(4, 7).Value = 0.3193946480751038
$ws.Cells.Item(4, 8).Value = 0.6879740953445435
$ws.Cells.Item(4, 9).Value = 0.4928436875343323

$ws.Cells.Item(5, 1).Value = "model_9_5_22"
$ws.Cells.Item(5, 2).Value = 0.8389620783597548
$ws.Cells.Item(5, 3).Value = 0.6954120463599007
$ws.Cells.Item(5, 4).Value = 0.2958535103341088
$ws.Cells.Item(5, 5).Value = 0.5351471719051999
$ws.Cells.Item(5, 6).Value = 0.1782215386629105
$ws.Cells.Item(5, 7).Value = 0.318206399679184
$ws.Cells.Item(5, 8).Value = 0.6886860728263855
$ws.Cells.Item(5, 9).Value = 0.4925496578216553

$ws.Cells.Item(6, 1).Value = "model_9_5_15"
$ws.Cells.Item(6, 2).Value = 0.8396524799798972
$ws.Cells.Item(6, 3).Value = 0.723133623107053
$ws.Cells.Item(6, 4).Value = 0.2769035587737952
$ws.Cells.Item(6, 5).Value = 0.5413858476635334
$ws.Cells.Item(6, 6).Value = 0.1774574518203735
$ws.Cells.Item(6, 7).Value = 0.2892453372478485
$ws.Cells.Item(6, 8).Value = 0.7072199583053589
$ws.Cells.Item(6, 9).Value = 0.4859392642974854

$ws.Cells.Item(7, 1).Value = "model_9_5_21"
$ws.Cells.Item(7, 2).Value = 0.8397871543282813
$ws.Cells.Item(7, 3).Value = 0.6975386231101776
$ws.Cells.Item(7, 4).Value = 0.2999968776081217
$ws.Cells.Item(7, 5).Value = 0.5380569541951716
$ws.Cells.Item(7, 6).Value = 0.1773084253072739
$ws.Cells.Item(7, 7).Value = 0.3159847855567932
$ws.Cells.Item(7, 8).Value = 0.6846336722373962
$ws.Cells.Item(7, 9).Value = 0.4894664883613586

$ws.Cells.Item(8, 1).Value = "model_9_5_23"
$ws.Cells.Item(8, 2).Value = 0.8397907632144309
$ws.Cells.Item(8, 3).Value = 0.6994013039949074
$ws.Cells.Item(8, 4).Value = 0.2973273813008819
$ws.Cells.Item(8, 5).Value = 0.5378696044415483
$ws.Cells.Item(8, 6).Value = 0.1773044019937515
$ws.Cells.Item(8, 7).Value = 0.3140388131141663
$ws.Cells.Item(8, 8).Value = 0.6872445344924927
$ws.Cells.Item(8, 9).Value = 0.4896650314331055

$ws.Cells.Item(9, 1).Value = "model_9_5_14"
$ws.Cells.Item(9, 2).Value = 0.84018657000726
$ws.Cells.Item(9, 3).Value = 0.7241550945945614
$ws.Cells.Item(9, 4).Value = 0.280158634629333
$ws.Cells.Item(9, 5).Value = 0.5433331927940408
$ws.Cells.Item(9, 6).Value = 0.1768663674592972
$ws.Cells.Item(9, 7).Value = 0.2881782352924347
$ws.Cells.Item(9, 8).Value = 0.7040363550186157
$ws.Cells.Item(9, 9).Value = 0.4838759005069733

$ws.Cells.Item(10, 1).Value = "model_9_5_24"
$ws.Cells.Item(10, 2).Value = 0.8402021772882751
$ws.Cells.Item(10, 3).Value = 0.7008309234805101
$ws.Cells.Item(10, 4).Value = 0.2985759682522142
$ws.Cells.Item(10, 5).Value = 0.539158178554034
$ws.Cells.Item(10, 6).Value = 0.1768490970134735
$ws.Cells.Item(10, 7).Value = 0.3125452399253845
$ws.Cells.Item(10, 8).Value = 0.6860233545303345
$ws.Cells.Item(10, 9).Value = 0.4882996380329132

$ws.Cells.Item(11, 1).Value = "model_9_5_12"
$ws.Cells.Item(11, 2).Value = 0.8403411661226758
$ws.Cells.Item(11, 3).Value = 0.7291997507779315
$ws.Cells.Item(11, 4).Value = 0.2791844685289004
$ws.Cells.Item(11, 5).Value = 0.5455427988970791
$ws.Cells.Item(11, 6).Value = 0.1766952723264694
$ws.Cells.Item(11, 7).Value = 0.2829080522060394
$ws.Cells.Item(11, 8).Value = 0.7049890756607056
$ws.Cells.Item(11, 9).Value = 0.4815346896648407

$ws.Cells.Item(12, 1).Value = "model_9_5_13"
$ws.Cells.Item(12, 2).Value = 0.8404094810515013
$ws.Cells.Item(12, 3).Value = 0.7222957669556245
$ws.Cells.Item(12, 4).Value = 0.2842859318134615
$ws.Cells.Item(12, 5).Value = 0.5441552721324268
$ws.Cells.Item(12, 6).Value = 0.176619678735733
$ws.Cells.Item(12, 7).Value = 0.2901206910610199
$ws.Cells.Item(12, 8).Value = 0.6999996900558472
$ws.Cells.Item(12, 9).Value = 0.4830048084259033

$ws.Cells.Item(13, 1).Value = "model_9_5_18"
$ws.Cells.Item(13, 2).Value = 0.8408079493895927
$ws.Cells.Item(13, 3).Value = 0.6970762149191156
$ws.Cells.Item(13, 4).Value = 0.3125221968343194
$ws.Cells.Item(13, 5).Value = 0.5432560322778508
$ws.Cells.Item(13, 6).Value = 0.1761786937713623
$ws.Cells.Item(13, 7).Value = 0.3164678514003754
$ws.Cells.Item(13, 8).Value = 0.6723833680152893
$ws.Cells.Item(13, 9).Value = 0.4839576780796051

$ws.Cells.Item(14, 1).Value = "model_9_5_16"
$ws.Cells.Item(14, 2).Value = 0.8412833057612754
$ws.Cells.Item(14, 3).Value = 0.7253215414830709
$ws.Cells.Item(14, 4).Value = 0.2864393442918494
$ws.Cells.Item(14, 5).Value = 0.5466700813750059
$ws.Cells.Item(14, 6).Value = 0.1756526082754135
$ws.Cells.Item(14, 7).Value = 0.2869596183300018
$ws.Cells.Item(14, 8).Value = 0.697893500328064
$ws.Cells.Item(14, 9).Value = 0.4803401827812195

$ws.Cells.Item(15, 1).Value = "model_9_5_11"
$ws.Cells.Item(15, 2).Value = 0.8413998442578267
$ws.Cells.Item(15, 3).Value = 0.7327726077002206
$ws.Cells.Item(15, 4).Value = 0.2833066186498167
$ws.Cells.Item(15, 5).Value = 0.5491985282980372
$ws.Cells.Item(15, 6).Value = 0.1755236238241196
$ws.Cells.Item(15, 7).Value = 0.2791754007339478
$ws.Cells.Item(15, 8).Value = 0.7009574770927429
$ws.Cells.Item(15, 9).Value = 0.4776611328125

$ws.Cells.Item(16, 1).Value = "model_9_5_17"
$ws.Cells.Item(16, 2).Value = 0.8416908946987862
$ws.Cells.Item(16, 3).Value = 0.7005604888362762
$ws.Cells.Item(16, 4).Value = 0.3164590446437542
$ws.Cells.Item(16, 5).Value = 0.5467849680312068
$ws.Cells.Item(16, 6).Value = 0.1752015203237534
$ws.Cells.Item(16, 7).Value = 0.3128277659416199
$ws.Cells.Item(16, 8).Value = 0.6685329675674438
$ws.Cells.Item(16, 9).Value = 0.4802184998989105

$ws.Cells.Item(17, 1).Value = "model_9_5_10"
$ws.Cells.Item(17, 2).Value = 0.8417154317231208
$ws.Cells.Item(17, 3).Value = 0.7368233970626197
$ws.Cells.Item(17, 4).Value = 0.2822531735872216
$ws.Cells.Item(17, 5).Value = 0.5508555080507106
$ws.Cells.Item(17, 6).Value = 0.1751743704080582
$ws.Cells.Item(17, 7).Value = 0.2749435007572174
$ws.Cells.Item(17, 8).Value = 0.7019877433776855
$ws.Cells.Item(17, 9).Value = 0.4759053587913513

$ws.Cells.Item(18, 1).Value = "model_9_5_9"
$ws.Cells.Item(18, 2).Value = 0.8471968432362729
$ws.Cells.Item(18, 3).Value = 0.7583141331038847
$ws.Cells.Item(18, 4).Value = 0.3061533238966727
$ws.Cells.Item(18, 5).Value = 0.5724547952079551
$ws.Cells.Item(18, 6).Value = 0.1691080778837204
$ws.Cells.Item(18, 7).Value = 0.2524918913841248
$ws.Cells.Item(18, 8).Value = 0.6786123514175415
$ws.Cells.Item(18, 9).Value = 0.4530191719532013

$ws.Cells.Item(19, 1).Value = "model_9_5_8"
$ws.Cells.Item(19, 2).Value = 0.8481195161235497
$ws.Cells.Item(19, 3).Value = 0.7603687204651407
$ws.Cells.Item(19, 4).Value = 0.3113660863907747
$ws.Cells.Item(19, 5).Value = 0.5757914422429034
$ws.Cells.Item(19, 6).Value = 0.1680869311094284
$ws.Cells.Item(19, 7).Value = 0.2503454685211182
$ws.Cells.Item(19, 8).Value = 0.6735141277313232
$ws.Cells.Item(19, 9).Value = 0.449483722448349

$ws.Cells.Item(20, 1).Value = "model_9_5_1"
$ws.Cells.Item(20, 2).Value = 0.8515598318919105
$ws.Cells.Item(20, 3).Value = 0.7706159695191763
$ws.Cells.Item(20, 4).Value = 0.5951439863642547
$ws.Cells.Item(20, 5).Value = 0.7044063754386047
$ws.Cells.Item(20, 6).Value = 0.1642795354127884
$ws.Cells.Item(20, 7).Value = 0.2396400570869446
$ws.Cells.Item(20, 8).Value = 0.3959668874740601
$ws.Cells.Item(20, 9).Value = 0.3132056593894958

$ws.Cells.Item(21, 1).Value = "model_9_5_5"
$ws.Cells.Item(21, 2).Value = 0.8543389337770623
$ws.Cells.Item(21, 3).Value = 0.828514870528325
$ws.Cells.Item(21, 4).Value = 0.3331380681400345
$ws.Cells.Item(21, 5).Value = 0.6208194118344839
$ws.Cells.Item(21, 6).Value = 0.1612038910388947
$ws.Cells.Item(21, 7).Value = 0.1791524291038513
$ws.Cells.Item(21, 8).Value = 0.6522201299667358
$ws.Cells.Item(21, 9).Value = 0.4017728865146637

$ws.Cells.Item(22, 1).Value = "model_9_5_7"
$ws.Cells.Item(22, 2).Value = 0.8571195319836498
$ws.Cells.Item(22, 3).Value = 0.8022401116732644
$ws.Cells.Item(22, 4).Value = 0.3602565608365512
$ws.Cells.Item(22, 5).Value = 0.6188845034755199
$ws.Cells.Item(22, 6).Value = 0.1581265777349472
$ws.Cells.Item(22, 7).Value = 0.2066019326448441
$ws.Cells.Item(22, 8).Value = 0.6256970763206482
$ws.Cells.Item(22, 9).Value = 0.4038230776786804

$ws.Cells.Item(23, 1).Value = "model_9_5_6"
$ws.Cells.Item(23, 2).Value = 0.8592788059116356
$ws.Cells.Item(23, 3).Value = 0.8465359501461609
$ws.Cells.Item(23, 4).Value = 0.3441920021586388
$ws.Cells.Item(23, 5).Value = 0.6350280415960599
$ws.Cells.Item(23, 6).Value = 0.155736893415451
$ws.Cells.Item(23, 7).Value = 0.1603256016969681
$ws.Cells.Item(23, 8).Value = 0.6414089202880859
$ws.Cells.Item(23, 9).Value = 0.3867177069187164

$ws.Cells.Item(24, 1).Value = "model_9_5_4"
$ws.Cells.Item(24, 2).Value = 0.8631508616369127
$ws.Cells.Item(24, 3).Value = 0.7789175909102253
$ws.Cells.Item(24, 4).Value = 0.5061207616409631
$ws.Cells.Item(24, 5).Value = 0.670070168334539
$ws.Cells.Item(24, 6).Value = 0.1514516770839691
$ws.Cells.Item(24, 7).Value = 0.2309672236442566
$ws.Cells.Item(24, 8).Value = 0.4830355048179626
$ws.Cells.Item(24, 9).Value = 0.3495876789093018

$ws.Cells.Item(25, 1).Value = "model_9_5_3"
$ws.Cells.Item(25, 2).Value = 0.8676089923079301
$ws.Cells.Item(25, 3).Value = 0.7669675457061988
$ws.Cells.Item(25, 4).Value = 0.5793066385925987
$ws.Cells.Item(25, 5).Value = 0.695622792787266
$ws.Cells.Item(25, 6).Value = 0.1465178281068802
$ws.Cells.Item(25, 7).Value = 0.2434515953063965
$ws.Cells.Item(25, 8).Value = 0.4114565253257751
$ws.Cells.Item(25, 9).Value = 0.3225125968456268

$ws.Cells.Item(26, 1).Value = "model_9_5_2"
$ws.Cells.Item(26, 2).Value = 0.8772586786608819
$ws.Cells.Item(26, 3).Value = 0.7819118475000837
$ws.Cells.Item(26, 4).Value = 0.6836700712254017
$ws.Cells.Item(26, 5).Value = 0.748756284738944
$ws.Cells.Item(26, 6).Value = 0.1358384788036346
$ws.Cells.Item(26, 7).Value = 0.2278390973806381
$ws.Cells.Item(26, 8).Value = 0.3093844950199127
$ws.Cells.Item(26, 9).Value = 0.2662133276462555

